# Auto-generated edit script: updates crypto price/volume table
# to match the target commit (cryptos list refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/string updates (Coin name, Link, Volume(1h)) ---
$plainUpdates = @(
  @{Cell='E2'; Value='  +0.36%  '},
  @{Cell='E3'; Value='  -0.02%  '},
  @{Cell='E4'; Value='  +0.00%  '},
  @{Cell='E5'; Value='  +0.10%  '},
  @{Cell='E6'; Value='  -0.67%  '},
  @{Cell='E7'; Value='  -0.04%  '},
  @{Cell='E8'; Value='  -0.34%  '},
  @{Cell='E9'; Value='  +1.69%  '},
  @{Cell='E10'; Value='  -0.40%  '},
  @{Cell='E11'; Value='  +6.39%  '},
  @{Cell='E12'; Value='  +3.65%  '},
  @{Cell='E13'; Value='  +0.01%  '},
  @{Cell='E14'; Value='  -1.91%  '},
  @{Cell='E15'; Value='  +1.67%  '},
  @{Cell='E16'; Value='  +0.16%  '},
  @{Cell='E17'; Value='  +0.60%  '},
  @{Cell='E18'; Value='  -0.77%  '},
  @{Cell='E19'; Value='  +0.88%  '},
  @{Cell='E20'; Value='  +0.57%  '},
  @{Cell='E21'; Value='  +2.65%  '},
  @{Cell='E22'; Value='  -0.24%  '},
  @{Cell='E23'; Value='  -1.83%  '},
  @{Cell='E24'; Value='  +0.02%  '},
  @{Cell='E25'; Value='  +4.55%  '},
  @{Cell='E26'; Value='  -0.49%  '},
  @{Cell='E27'; Value='  -0.32%  '},
  @{Cell='E28'; Value='  +0.03%  '},
  @{Cell='E29'; Value='  -0.08%  '},
  @{Cell='E30'; Value='  +0.46%  '},
  @{Cell='B31'; Value='Filecoin'},
  @{Cell='C31'; Value='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'},
  @{Cell='E31'; Value='  +1.76%  '},
  @{Cell='B32'; Value='PancakeSwap'},
  @{Cell='C32'; Value='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'},
  @{Cell='E32'; Value='  -0.99%  '},
  @{Cell='E33'; Value='  +0.72%  '},
  @{Cell='E34'; Value='  -1.19%  '},
  @{Cell='E35'; Value='  -1.28%  '},
  @{Cell='E36'; Value='  -0.33%  '},
  @{Cell='E37'; Value='  +0.39%  '},
  @{Cell='E38'; Value='  +2.35%  '},
  @{Cell='E39'; Value='  +2.52%  '},
  @{Cell='E40'; Value='  -0.05%  '},
  @{Cell='E41'; Value='  +1.12%  '},
  @{Cell='E42'; Value='  -3.61%  '},
  @{Cell='E43'; Value='  -0.34%  '},
  @{Cell='E44'; Value='  +3.29%  '},
  @{Cell='E45'; Value='  +0.32%  '},
  @{Cell='E46'; Value='  -0.09%  '},
  @{Cell='E47'; Value='  +0.26%  '},
  @{Cell='E48'; Value='  -0.38%  '},
  @{Cell='E49'; Value='  +2.54%  '},
  @{Cell='E50'; Value='  +1.35%  '},
  @{Cell='E51'; Value='  -0.92%  '}
)

foreach ($u in $plainUpdates) {
  $ws.Range($u.Cell).Value = $u.Value
}

# --- Price (column D) updates: force text storage so values like
#     "29.897.02" or "25.61" are not reinterpreted as numbers/dates ---
$priceUpdates = @(
  @{Cell='D2'; Value='29.897.02'},
  @{Cell='D3'; Value='1.888.01'},
  @{Cell='D5'; Value='0.7721'},
  @{Cell='D6'; Value='242.81'},
  @{Cell='D8'; Value='0.3115'},
  @{Cell='D9'; Value='25.61'},
  @{Cell='D10'; Value='0.07188'},
  @{Cell='D11'; Value='0.08599'},
  @{Cell='D12'; Value='1.960.14'},
  @{Cell='D13'; Value='0.7645'},
  @{Cell='D14'; Value='5.375'},
  @{Cell='D15'; Value='93.79'},
  @{Cell='D16'; Value='6.188'},
  @{Cell='D17'; Value='29.965.42'},
  @{Cell='D19'; Value='244.51'},
  @{Cell='D20'; Value='0.000007815'},
  @{Cell='D21'; Value='2.201.97'},
  @{Cell='D22'; Value='0.9986'},
  @{Cell='D23'; Value='7.997'},
  @{Cell='D25'; Value='0.1651'},
  @{Cell='D26'; Value='9.374'},
  @{Cell='D27'; Value='161.99'},
  @{Cell='D28'; Value='18.74'},
  @{Cell='D29'; Value='2.032'},
  @{Cell='D30'; Value='1.443'},
  @{Cell='D31'; Value='4.532'},
  @{Cell='D32'; Value='1.533'},
  @{Cell='D33'; Value='4.104'},
  @{Cell='D34'; Value='0.05429'},
  @{Cell='D35'; Value='1.241'},
  @{Cell='D36'; Value='0.7459'},
  @{Cell='D38'; Value='2.695'},
  @{Cell='D39'; Value='0.01968'},
  @{Cell='D41'; Value='0.4465'},
  @{Cell='D42'; Value='1.109.29'},
  @{Cell='D43'; Value='73.37'},
  @{Cell='D45'; Value='0.8508'},
  @{Cell='D48'; Value='1.872'},
  @{Cell='D49'; Value='7.626'},
  @{Cell='D50'; Value='2.076.44'},
  @{Cell='D51'; Value='2.985'}
)

foreach ($u in $priceUpdates) {
  $c = $ws.Range($u.Cell)
  $c.NumberFormat = "@"
  $c.Value = $u.Value
  $c.Style = "Normal"
}
